$d = $word.ActiveDocument

# This document contains two occurrences of the literal folder name "Assets"
# (inside Python string literals shown in a code listing) that must become
# "assets". A plain text Find/Replace would make the host's run-merging
# collapse the edited text into a single run together with whatever
# formatting-identical runs happen to sit next to it, which does not match
# the very specific run layout produced by the original (interactive) edit:
# the author selected just the capital "A" and retyped a lowercase "a",
# which made Word split the run into three pieces - the text before the
# edit point, the freshly-typed "a", and the text after the edit point -
# while leaving any other already-separate runs alone.
#
# To reproduce that exact run layout we change the single "A" character
# using a temporary, distinguishing font color so the freshly created runs
# cannot silently re-merge with neighboring runs that already share the
# "../Assets/..." run's formatting, then restore the real color
# (067D17 == 1539334 in the BGR/wdColor form Word uses) once every run
# involved has been created/edited.
$greenBgr = 1539334   # 0x177D06 -> matches <w:color w:val="067D17"/>
$tempBgr1 = 1539333
$tempBgr2 = 1234567

function Fix-AssetsOccurrence($AIndex, $RestLength, $RestText, $ProtectNext, $NextLength) {
    $nextStart = $AIndex + 1 + $RestLength

    if ($ProtectNext) {
        # Temporarily detune the following sibling run's color so it cannot
        # be folded back together with the run we are about to rebuild.
        $rNext = $d.Range($nextStart, $nextStart + $NextLength)
        $rNext.Font.Color = $tempBgr2
    }

    # Split "a" off from the rest of the word by giving it a distinguishing
    # color while retyping it.
    $rA = $d.Range($AIndex, $AIndex + 1)
    $rA.Font.Color = $tempBgr1
    $rA.Text = "a"

    # Rewrite the remainder ("ssets...") via a throwaway placeholder first -
    # this keeps it from being silently reabsorbed into the run to its left,
    # and marks it as a genuinely new run (so it ends up without any rsid
    # attributes, matching a freshly authored edit).
    $restStart = $AIndex + 1
    $restEnd = $restStart + $RestLength
    $rRest = $d.Range($restStart, $restEnd)
    $placeholder = ""
    for ($i = 0; $i -lt $RestLength; $i++) {
        $placeholder += "z"
    }
    $rRest.Font.Color = $tempBgr1
    $rRest.Text = $placeholder

    $rRest2 = $d.Range($restStart, $restEnd)
    $rRest2.Font.Color = $greenBgr
    $rRest2.Text = $RestText

    # Restore the real color on the new "a" run.
    $rA2 = $d.Range($AIndex, $AIndex + 1)
    $rA2.Font.Color = $greenBgr

    if ($ProtectNext) {
        # Restore the sibling run's original color.
        $rNext2 = $d.Range($nextStart, $nextStart + $NextLength)
        $rNext2.Font.Color = $greenBgr
    }
}

# --- Occurrence 1: ../Assets/    (already followed by a separate "fighter.png"" run) ---
$full = $d.Content.Text
$idx1 = $full.IndexOf("../Assets/")
$aIndex1 = $idx1 + 3   # "../" is 3 characters before the "A"
Fix-AssetsOccurrence $aIndex1 6 "ssets/" $true 12

# --- Occurrence 2: "../Assets/pew.wav"   (the whole quoted string is one run) ---
$full2 = $d.Content.Text
$idx2 = $full2.IndexOf('"../Assets/pew.wav"')
$aIndex2 = $idx2 + 4   # '"../' is 4 characters before the "A"
Fix-AssetsOccurrence $aIndex2 14 'ssets/pew.wav"' $false 0

Write-Output "done"
